# Scrape data project update
# Rework invoice layout: move the item table down to start at row 5,
# rename/re-letter the headers, add SubTotal/GST/Total columns (G/H/I)
# to the header row, and shift the item rows + footer rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture existing item-row & footer values before we overwrite anything ---
# NOTE: use .Value2 for reads -- .Value as a bare read returns a descriptor
# stub in this host; .Value2 returns the actual scalar.
$items = @()
for ($r = 2; $r -le 8; $r++) {
    $items += ,@(
        $ws.Cells.Item($r, 3).Value2,  # C - ID/ItemNo
        $ws.Cells.Item($r, 4).Value2,  # D - Description
        $ws.Cells.Item($r, 5).Value2,  # E - Quantity
        $ws.Cells.Item($r, 6).Value2,  # F - Price
        $ws.Cells.Item($r, 7).Value2   # G - Total
    )
}

$footer = @()
for ($r = 9; $r -le 11; $r++) {
    $footer += ,@(
        $ws.Cells.Item($r, 6).Value2,  # F - label
        $ws.Cells.Item($r, 7).Value2   # G - value
    )
}
$subTotalValue = $footer[0][1]
$gstValue      = $footer[1][1]
$totalValue    = $footer[2][1]

# --- clear out the old table area (rows 1-11, columns C-G) ---
$ws.Range("C1:G11").Clear() | Out-Null

# --- new header row at row 5 (columns C:I) ---
$ws.Cells.Item(5, 3).Value = "ItemNo"
$ws.Cells.Item(5, 4).Value = "Description"
$ws.Cells.Item(5, 5).Value = "Quantity"
$ws.Cells.Item(5, 6).Value = "Price"
$ws.Cells.Item(5, 7).Value = "SubTotal"
$ws.Cells.Item(5, 8).Value = "GST"
$ws.Cells.Item(5, 9).Value = "Total"

# --- rewrite the item rows shifted down by 4 (old row N -> new row N+4) ---
for ($i = 0; $i -lt $items.Count; $i++) {
    $oldRow = 2 + $i
    $newRow = $oldRow + 4
    $row = $items[$i]
    $ws.Cells.Item($newRow, 3).Value = $row[0]
    $ws.Cells.Item($newRow, 4).Value = $row[1]
    $ws.Cells.Item($newRow, 5).Value = $row[2]
    $ws.Cells.Item($newRow, 6).Value = $row[3]
    $ws.Cells.Item($newRow, 7).Value = $row[4]
}

# first item row (row 6) also carries the GST and Total figures
$ws.Cells.Item(6, 8).Value = $gstValue
$ws.Cells.Item(6, 9).Value = $totalValue

# --- rewrite footer rows shifted down by 4 (old row N -> new row N+4) ---
$ws.Cells.Item(13, 6).Value = "Sub Total "
$ws.Cells.Item(13, 7).Value = $subTotalValue

$ws.Cells.Item(14, 6).Value = "GST 8% "
$ws.Cells.Item(14, 7).Value = $gstValue

$ws.Cells.Item(15, 6).Value = "Total "
$ws.Cells.Item(15, 7).Value = $totalValue
